$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table data (header stays the same; body is fully replaced/reordered)
$data = @(
    @("Payton Pritchard",     "PG,SG",    "Boston Celtics"),
    @("Donovan Mitchell",     "PG,SG",    "Cleveland Cavaliers"),
    @("Josh Hart",            "SG,SF,PF", "New York Knicks"),
    @("Malik Beasley",        "SG,SF",    "Detroit Pistons"),
    @("Michael Porter Jr.",   "SF,PF",    "Denver Nuggets"),
    @("P.J. Washington",      "PF",       "Dallas Mavericks"),
    @("Alexandre Sarr",       "PF,C",     "Washington Wizards"),
    @("De'Andre Hunter",      "SF,PF",    "Atlanta Hawks"),
    @("Domantas Sabonis",     "C",        "Sacramento Kings"),
    @("Alperen Sengün",       "C",        "Houston Rockets"),
    @("Dyson Daniels",        "PG,SG,SF", "Atlanta Hawks"),
    @("Victor Wembanyama",    "C",        "San Antonio Spurs"),
    @("Cam Thomas",           "SG,SF",    "Brooklyn Nets"),
    @("Shaedon Sharpe",       "SG,SF",    "Portland Trail Blazers"),
    @("Deandre Ayton",        "C",        "Portland Trail Blazers"),
    @("Kristaps Porzingis",   "PF,C",     "Boston Celtics"),
    @("Tari Eason",           "SF,PF",    "Houston Rockets"),
    @("Bradley Beal",         "PG,SG,SF", "Phoenix Suns")
)

$rowCount = $data.Length
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
